$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.338.50'
$ws.Range('E2').Value = '  +0.38%  '

$ws.Range('D3').Value = '1.861.24'
$ws.Range('E3').Value = '  +0.10%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.17%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7025'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.12%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '238.22'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.28%  '

$ws.Range('E7').Value = '  +0.10%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07881'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.12%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3050'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.01%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.78'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +6.36%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08170'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.26%  '

$ws.Range('D12').Value = '1.882.03'
$ws.Range('E12').Value = '  +1.54%  '

$ws.Range('E13').Value = '  +0.93%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7145'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.62%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.34'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.25%  '

$ws.Range('D16').Value = '29.431.75'
$ws.Range('E16').Value = '  +0.68%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.825'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.83%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000007789'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.10%  '

$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '239.02'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.16%  '

$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.22'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.12%  '

$ws.Range('D21').Value = '2.157.62'
$ws.Range('E21').Value = '  +2.96%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.15%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.001'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.19%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.537'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.19%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '162.74'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.60%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.897'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -1.20%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.1421'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -3.21%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.07'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.23%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.904'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -5.45%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.371'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -4.51%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.475'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.74%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.305'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.62%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.046'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -0.12%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05172'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.95%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.180'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.85%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7077'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.10%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.002'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.04%  '

$ws.Range('E38').Value = '  +0.24%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01841'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.27%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.695'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.87%  '

$ws.Range('D41').Value = '1.171.45'
$ws.Range('E41').Value = '  +2.74%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.9207'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.69%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.040'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.97%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '71.79'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +1.38%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4250'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.68%  '

$ws.Range('E46').Value = '  +0.13%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.80'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.61%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.5352'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.54%  '

$ws.Range('E49').Value = '  -2.38%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.166'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.26%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.987'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.08%  '
